$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a few existing score values for 2025-02-04 (rows 14 & 15) ---
$ws.Range("E14").Value = 10
$ws.Range("G14").Value = 9.9963696766329875
$ws.Range("K14").Value = 5.838972561904181
$ws.Range("M14").Value = 8.8576098012089517
$ws.Range("P14").Value = 51.186996347884147

$ws.Range("E15").Value = 8.6531568789633315
$ws.Range("K15").Value = 6.8340438912736658
$ws.Range("P15").Value = 45.010292575777243

# --- Narrow column A a touch ---
$ws.Columns("A").ColumnWidth = 22

# --- Append the new daily rows (2025-02-05) ---
# The date column has to land as plain text (matching the existing "Date"
# column cells) rather than be auto-recognised as a real date, so build the
# text via a formula first and then flatten it to a static value.
$ws.Range("A18").Formula = "=""2025-02-05"""
$ws.Range("A19").Formula = "=""2025-02-05"""
$ws.Range("A20").Formula = "=""2025-02-05"""
$ws.Range("A21").Formula = "=""2025-02-05"""
$ws.Range("A18:A21").Copy() | Out-Null
$ws.Range("A18:A21").PasteSpecial(-4163) | Out-Null

$ws.Range("B18").Value = "abs_activity"
$ws.Range("B19").Value = "rel_activity"
$ws.Range("B20").Value = "abs_sleep"
$ws.Range("B21").Value = "rel_sleep"

$ws.Range("C18").Value = 9.1538272296850707
$ws.Range("D18").Value = 5.0079195451244303
$ws.Range("E18").Value = 9.0995955873121481
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9.3928467163721248
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 8.8410306119502025
$ws.Range("K18").Value = 9.7067890020044469
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = 8.6442053003989052
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 45.997263835772699
$ws.Range("Q18").Value = 43.848950157074633

$ws.Range("C19").Value = 7.9418103829017426
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 6.2406947890818856
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 10
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 9.3688320075931593
$ws.Range("K19").Value = 10
$ws.Range("L19").Value = 10
$ws.Range("M19").Value = 5.7342619070832113
$ws.Range("N19").Value = 5
$ws.Range("O19").Value = 5
$ws.Range("P19").Value = 34.916767079066837
$ws.Range("Q19").Value = 49.368832007593163

$ws.Range("C20").Value = 9.5666666666666664
$ws.Range("D20").Value = 9.4666666666666668
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 8.1
$ws.Range("G20").Value = 8.6666666666666661
$ws.Range("H20").Value = 1.8666666666666669
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 9.6
$ws.Range("K20").Value = 7.6000000000000014
$ws.Range("L20").Value = 10
$ws.Range("M20").Value = 10
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 45.833333333333343
$ws.Range("Q20").Value = 39.033333333333331

$ws.Range("C21").Value = 7.8773084177742678
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 8.0124654398816197
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 7
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 7
$ws.Range("M21").Value = 7
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 22.889773857655889
$ws.Range("Q21").Value = 14

# --- Match the selection Excel leaves behind after typing the new block ---
$ws.Range("A22:Q25").Select() | Out-Null
